# Refresh the cryptos snapshot: updated prices / 1h volume % for most rows,
# plus two pairs of rows that swapped rank order (ImmutableX/Stellar and
# FraxShare/ARBITRUM).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) mixes "thousands.decimal" style numbers with plain
# decimals and must stay plain text (as in the source data), so force a
# text format while assigning values, then restore the default style.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "37.668.31"
$ws.Range("E2").Value = "  +6.19%  "
$ws.Range("D3").Value = "2.057.56"
$ws.Range("E3").Value = "  +3.35%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "252.86"
$ws.Range("E5").Value = "  +4.96%  "
$ws.Range("E6").Value = "  +2.87%  "
$ws.Range("D7").Value = "66.69"
$ws.Range("E7").Value = "  +18.42%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  +6.12%  "
$ws.Range("D10").Value = "59.47"
$ws.Range("E10").Value = "  +0.65%  "
$ws.Range("E11").Value = "  +4.56%  "
$ws.Range("E12").Value = "  +1.03%  "
$ws.Range("D13").Value = "15.32"
$ws.Range("E13").Value = "  +6.92%  "
$ws.Range("E14").Value = "  +2.20%  "
$ws.Range("D15").Value = "2.357.09"
$ws.Range("E15").Value = "  +3.32%  "
$ws.Range("E16").Value = "  +7.51%  "
$ws.Range("D17").Value = "20.81"
$ws.Range("E17").Value = "  +22.38%  "
$ws.Range("D18").Value = "2.050.25"
$ws.Range("E18").Value = "  +3.04%  "
$ws.Range("D19").Value = "37.563.71"
$ws.Range("E19").Value = "  +5.91%  "
$ws.Range("D20").Value = "73.35"
$ws.Range("E20").Value = "  +5.09%  "
$ws.Range("D21").Value = "0.0₃0880"
$ws.Range("E21").Value = "  +5.56%  "
$ws.Range("E22").Value = "  +7.18%  "
$ws.Range("D23").Value = "238.00"
$ws.Range("E23").Value = "  +2.84%  "
$ws.Range("D24").Value = "2.78"
$ws.Range("E24").Value = "  +23.17%  "
$ws.Range("E25").Value = "  -0.11%  "
$ws.Range("E26").Value = "  +4.31%  "
$ws.Range("D27").Value = "9.63"
$ws.Range("E27").Value = "  +5.74%  "
$ws.Range("D28").Value = "165.60"
$ws.Range("E28").Value = "  +1.60%  "
$ws.Range("D29").Value = "19.97"
$ws.Range("E29").Value = "  +2.78%  "
$ws.Range("D30").Value = "5.23"
$ws.Range("E30").Value = "  +9.97%  "
$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").Value = "0.122"
$ws.Range("E31").Value = "  +3.10%  "
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "1.23"
$ws.Range("E32").Value = "  +8.07%  "
$ws.Range("D33").Value = "0.112"
$ws.Range("E33").Value = "  +24.66%  "
$ws.Range("D34").Value = "4.76"
$ws.Range("E34").Value = "  +12.10%  "
$ws.Range("E35").Value = "  +5.60%  "
$ws.Range("E36").Value = "  +10.17%  "
$ws.Range("E37").Value = "  -0.14%  "
$ws.Range("D38").Value = "6.05"
$ws.Range("E38").Value = "  +24.66%  "
$ws.Range("E39").Value = "  +0.84%  "
$ws.Range("E40").Value = "  +17.20%  "
$ws.Range("E41").Value = "  +5.46%  "
$ws.Range("E42").Value = "  +4.40%  "
$ws.Range("D43").Value = "0.0220"
$ws.Range("E43").Value = "  +5.49%  "
$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D44").Value = "1.14"
$ws.Range("E44").Value = "  +6.11%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").Value = "8.16"
$ws.Range("E45").Value = "  +10.19%  "
$ws.Range("D46").Value = "17.17"
$ws.Range("E46").Value = "  +11.39%  "
$ws.Range("E47").Value = "  +20.07%  "
$ws.Range("D48").Value = "95.83"
$ws.Range("E48").Value = "  +6.39%  "
$ws.Range("D49").Value = "1.432.06"
$ws.Range("E49").Value = "  +4.76%  "
$ws.Range("D50").Value = "2.95"
$ws.Range("E50").Value = "  +2.07%  "
$ws.Range("D51").Value = "47.70"
$ws.Range("E51").Value = "  +5.35%  "

$ws.Range("D2:D51").Style = "Normal"
